$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.984.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.740.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9976"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4979"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07270"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.060"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9968"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.896"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.731.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.838"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06394"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9977"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.020.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.058"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.926.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.095"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.053"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.575"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.410"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05934"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02190"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.426"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.762"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6029"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9986"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.114"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.577"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5638"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "119.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.847"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06670"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.100"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.43%  "
